$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.1
$ws.Range("L2").Value = 4.5
$ws.Range("N2").Value = 21
$ws.Range("S2").Value = 1.22
$ws.Range("T2").Value = 4
$ws.Range("X2").Value = 10
$ws.Range("AB2").Value = 17
$ws.Range("AH2").Value = 21
$ws.Range("AL2").Value = 34
$ws.Range("AO2").Value = 8
$ws.Range("AQ2").Value = 21
$ws.Range("AR2").Value = 34
$ws.Range("AT2").Value = 4
$ws.Range("AX2").Value = 7
$ws.Range("BC2").Value = 126
